$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2750
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2325
$ws.Range("N40").Value = -3350
$ws.Range("H41").Value = 420.6842
$ws.Range("I41").Value = 378.08334
$ws.Range("K41").Value = 378.08334
$ws.Range("M41").Value = 61.91665999999998
$ws.Range("H116").Value = 5142
$ws.Range("J116").Value = 4415.6665
$ws.Range("L116").Value = 4415.6665
$ws.Range("N116").Value = -11299.6665
$ws.Range("H132").Value = 2209.121
$ws.Range("I132").Value = 2165.7273
$ws.Range("J132").Value = 2295.9092
$ws.Range("K132").Value = 6497.1819
$ws.Range("L132").Value = 6887.7276
$ws.Range("M132").Value = -3967.1819
$ws.Range("N132").Value = -11947.7276
$ws.Range("H135").Value = 163.66667
$ws.Range("I135").Value = 163.66667
$ws.Range("K135").Value = 1473.00003
$ws.Range("M135").Value = 1061.99997
$ws.Range("H137").Value = 2275
$ws.Range("I137").Value = 1975.25
$ws.Range("J137").Value = 2574.75
$ws.Range("K137").Value = 5925.75
$ws.Range("L137").Value = 7724.25
$ws.Range("M137").Value = -3375.75
$ws.Range("N137").Value = -12824.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 49897
$ws.Range("J130").Value = 49897
$ws.Range("L130").Value = 49897
$ws.Range("N130").Value = -59937
$ws.Range("H132").Value = 1301.8
$ws.Range("I132").Value = 1301.8
$ws.Range("K132").Value = 3905.4
$ws.Range("M132").Value = -1375.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1312.5
$ws.Range("I134").Value = 1101.5625
$ws.Range("K134").Value = 3304.6875
$ws.Range("M134").Value = -769.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2297.8572
$ws.Range("I31").Value = 2218
$ws.Range("J31").Value = 2404.3333
$ws.Range("K31").Value = 2218
$ws.Range("L31").Value = 2404.3333
$ws.Range("M31").Value = -1923
$ws.Range("N31").Value = -2994.3333
$ws.Range("H34").Value = 2297.8572
$ws.Range("I34").Value = 2218
$ws.Range("J34").Value = 2404.3333
$ws.Range("K34").Value = 2218
$ws.Range("L34").Value = 2404.3333
$ws.Range("M34").Value = -2016
$ws.Range("N34").Value = -2808.3333
$ws.Range("H58").Value = 1495
$ws.Range("I58").Value = 1495
$ws.Range("K58").Value = 1495
$ws.Range("M58").Value = -1292
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("N91").Value = 0
$ws.Range("H105").Value = 2150.5
$ws.Range("I105").Value = 2060.8
$ws.Range("K105").Value = 2060.8
$ws.Range("M105").Value = -313.8000000000002
$ws.Range("H122").Value = 899.61536
$ws.Range("J122").Value = 910.7143
$ws.Range("L122").Value = 2732.1429
$ws.Range("N122").Value = -7632.1429
$ws.Range("H134").Value = 2233.238
$ws.Range("I134").Value = 1975
$ws.Range("K134").Value = 5925
$ws.Range("M134").Value = -3390
$ws.Range("H136").Value = 1495
$ws.Range("I136").Value = 1495
$ws.Range("K136").Value = 4485
$ws.Range("M136").Value = -1935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 987
$ws.Range("I36").Value = 224
$ws.Range("J36").Value = 1750
$ws.Range("K36").Value = 672
$ws.Range("L36").Value = 5250
$ws.Range("M36").Value = -503
$ws.Range("N36").Value = -5588

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 50000
$ws.Range("J34").Value = 50000
$ws.Range("L34").Value = 50000
$ws.Range("N34").Value = -50536
$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50630
$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52184
$ws.Range("H97").Value = 770.53845
$ws.Range("I97").Value = 475.55554
$ws.Range("K97").Value = 475.55554
$ws.Range("M97").Value = 20.44445999999999
$ws.Range("H113").Value = 1879.3334
$ws.Range("I113").Value = 1774.8572
$ws.Range("K113").Value = 1774.8572
$ws.Range("M113").Value = 395.1428000000001
$ws.Range("H132").Value = 5916.6665
$ws.Range("I132").Value = 5916.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17749.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -15219.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4999.8335
$ws.Range("I40").Value = 5054.8887
$ws.Range("K40").Value = 5054.8887
$ws.Range("M40").Value = -4918.8887
$ws.Range("H46").Value = 1631.3667
$ws.Range("I46").Value = 1389.2858
$ws.Range("J46").Value = 1843.1875
$ws.Range("K46").Value = 1389.2858
$ws.Range("L46").Value = 1843.1875
$ws.Range("M46").Value = -1201.2858
$ws.Range("N46").Value = -2219.1875
$ws.Range("H68").Value = 2071.6667
$ws.Range("I68").Value = 2130
$ws.Range("J68").Value = 1998.75
$ws.Range("K68").Value = 2130
$ws.Range("L68").Value = 1998.75
$ws.Range("M68").Value = -1381
$ws.Range("N68").Value = -3496.75
$ws.Range("H71").Value = 2071.6667
$ws.Range("I71").Value = 2130
$ws.Range("J71").Value = 1998.75
$ws.Range("K71").Value = 10650
$ws.Range("L71").Value = 9993.75
$ws.Range("M71").Value = -6906
$ws.Range("N71").Value = -17481.75
$ws.Range("H136").Value = 4334.222
$ws.Range("I136").Value = 3858.2856
$ws.Range("K136").Value = 11574.8568
$ws.Range("M136").Value = -9024.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4666.3335
$ws.Range("J62").Value = 5249.5
$ws.Range("L62").Value = 5249.5
$ws.Range("N62").Value = -6497.5
$ws.Range("H65").Value = 4666.3335
$ws.Range("J65").Value = 5249.5
$ws.Range("L65").Value = 26247.5
$ws.Range("N65").Value = -32487.5
$ws.Range("H80").Value = 24833
$ws.Range("J80").Value = 24833
$ws.Range("L80").Value = 24833
$ws.Range("N80").Value = -26829
$ws.Range("H81").Value = 10000
$ws.Range("J81").Value = 11000
$ws.Range("L81").Value = 22000
$ws.Range("N81").Value = -24122
$ws.Range("H83").Value = 24833
$ws.Range("J83").Value = 24833
$ws.Range("L83").Value = 74499
$ws.Range("N83").Value = -84483
$ws.Range("H84").Value = 10000
$ws.Range("J84").Value = 11000
$ws.Range("L84").Value = 110000
$ws.Range("N84").Value = -120608
$ws.Range("H100").Value = 5809863
$ws.Range("I100").Value = 8713071
$ws.Range("K100").Value = 17426142
$ws.Range("M100").Value = -17425601
$ws.Range("H107").Value = 482.33334
$ws.Range("I107").Value = 398.75
$ws.Range("K107").Value = 1196.25
$ws.Range("M107").Value = 723.75
$ws.Range("H113").Value = 601.75
$ws.Range("J113").Value = 825
$ws.Range("L113").Value = 2475
$ws.Range("N113").Value = -6815
$ws.Range("H122").Value = 1184.1428
$ws.Range("I122").Value = 1260.909
$ws.Range("K122").Value = 3782.727
$ws.Range("M122").Value = -1332.727
$ws.Range("H132").Value = 1965.04
$ws.Range("I132").Value = 1739.8572
$ws.Range("J132").Value = 3147.25
$ws.Range("K132").Value = 5219.571599999999
$ws.Range("L132").Value = 9441.75
$ws.Range("M132").Value = -2689.571599999999
$ws.Range("N132").Value = -14501.75
$ws.Range("H136").Value = 1950.5238
$ws.Range("I136").Value = 1464.1538
$ws.Range("J136").Value = 2740.875
$ws.Range("K136").Value = 4392.4614
$ws.Range("L136").Value = 8222.625
$ws.Range("M136").Value = -1842.4614
$ws.Range("N136").Value = -13322.625
